# Removed all the pT values less than 1 GeV.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every pT value (column B, rows 2-17) up by one GeV bin, so the
# row that used to hold pT=0 now holds pT=1, ..., up to pT=16.
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = $cell.Value2 + 1
}

# The old pT=16 row (row 18) is no longer needed - its bin is now
# represented by the shifted row 17. Clear its contents but keep the
# existing cell formatting in place.
$ws.Range("A18:J18").ClearContents()

# Update the active selection to match the new data block.
$ws.Range("A2:J17").Select()
